$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title ---
Replace-Text "Quasar Pulsations: Echoes of the Early Universe" "Understanding the Beauty of Chemistry: A Journey into the World of Matter"

# --- Author name (merges 3 runs "Dr" + "." + " Elara Jamil" into one) ---
Replace-Text "Dr. Elara Jamil" "Daniel Clark"

# --- Email (merges 3 runs "elarajamil@avatarastrophysics" + "." + "org" into one) ---
Replace-Text "elarajamil@avatarastrophysics.org" "dot"

# --- Body paragraph 1 ---
Replace-Text "In the vast cosmic tapestry, quasars, the luminous beacons of distant galaxies, captivate our imagination" "Chemistry, often perceived as a daunting subject, is in fact an enchanting realm that unravels the secrets of matter and its interactions"

Replace-Text " These brilliant celestial wonders, powered by supermassive black holes, emit immense amounts of energy across the electromagnetic spectrum" " It is a science that touches every aspect of our existence, from the air we breathe to the food we eat, and holds the key to unlocking countless mysteries of the universe"

Replace-Text " Their enigmatic pulsations, like rhythmic heartbeats of the cosmos, hold clues to unraveling the mysteries of the early universe" " In this exploration, we embark on a journey into the captivating world of chemistry, where we unravel the structure of matter, delve into the complexities of chemical reactions, and discover the profound impact chemistry has on our daily lives and the world around us"

Replace-Text "Peering into the annals of time, astronomers have discovered quasars that pulsate with remarkable regularity" "Chemistry is the study of matter and its properties, revealing the composition, structure, behavior, and transformations of substances"

Replace-Text " These pulsations, manifested as variations in brightness, offer a unique window into the dynamic interplay between the supermassive black hole and its surrounding accretion disk" " It explores the fundamental principles that govern the interactions between atoms and molecules, providing insights into the innermost workings of matter and its intricate relationships"

Replace-Text " As matter spirals inward, like a cosmic ballet, it releases phenomenal amounts of energy, giving rise to the rhythmic pulsations that enthrall astronomers" " Through chemical reactions, we witness the rearrangement of atoms, the formation and breaking of bonds, and the release or absorption of energy, uncovering the secrets of chemical reactivity and the driving forces behind chemical change"

Replace-Text "Moreover, the pulsations of quasars provide valuable insights into the birth and evolution of galaxies" "The world of chemistry is a tapestry of colors, textures, and aromas, each holding clues to the composition and properties of different substances"

Replace-Text " By analyzing the patterns and characteristics of these cosmic drumbeats, astrophysicists can probe the properties of the host galaxies, unveil the mysteries of black hole growth, and gain deeper insights into the distant epochs of cosmic history" " From the vibrant hues of transition metal complexes to the distinct odors of organic compounds, chemistry engages our senses, captivating our imaginations and stimulating our curiosity"

Replace-Text " Quasars serve as distant lighthouses, guiding us through the murky depths of time and space, illuminating the formative stages of the universe" " It challenges us to unravel the enigmas of matter, to understand the intricate dance of molecules, and to harness the power of chemistry to address global challenges and improve human lives"

# --- Summary paragraph ---
Replace-Text "The rhythmic pulsations of quasars, like celestial metronomes, offer a mesmerizing glimpse into the enigmatic depths of the early universe" "Our exploration of chemistry has illuminated the captivating world of matter and its interactions, unveiling the fundamental principles governing chemical reactions and the profound impact chemistry has on our daily lives"

Replace-Text " Their pulsations, echoing across cosmic eons, provide a unique probe to study the dynamic interplay between supermassive black holes and their accretion disks" " From the composition of substances to the transformations they undergo, chemistry provides a lens through which we gain insights into the innermost workings of the universe"

Replace-Text " Furthermore, these cosmic drumbeats hold valuable clues to unraveling the mysteries of galaxy formation and evolution" " It is a science that stimulates our senses, challenges our intellect, and holds the potential to solve some of the world's most pressing challenges"

Replace-Text " As we continue to decipher the intricate patterns of quasar pulsations, we unlock profound secrets about the birth and maturation of galaxies, enriching our understanding of the grand cosmic narrative" " As we continue to unravel the mysteries of chemistry, we embark on a journey of discovery, unlocking the secrets of matter and its transformative power"

# --- Apply Times New Roman font across the whole document (fixes "TimesNewToman" typo) ---
$r = $d.Range(0, $d.Content.End)
$r.Font.Name = "Times New Roman"
